# Fixed problem with calculating r from model
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Model results")

# Row 2: clear G2 and H2 (previously had r.model.h / r.model.f values)
$ws.Range("G2:H2").ClearContents()

# Row 3: update G3 (r.model.h) and H3 (r.model.f) with recalculated values,
# and clear the note in I3 that explained the old (incorrect) calculation
$ws.Range("G3").Value = 0.0232
$ws.Range("H3").Value = 0.0275
$ws.Range("I3").ClearContents()

# Row 4: clear G4 (r.model.h) and H4 (r.model.f) values
$ws.Range("G4").ClearContents()
$ws.Range("H4").ClearContents()

# Update the active selection to match the edited workbook
$ws.Range("H4").Select()
